# Enregistrement des articles dans le tableau pour les particuliers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Menus Crêpes" column entries that were mistakenly entered
# for the "particuliers" rows (Truffes, Meringues, Kig ha farz,
# Potée de pouldrezic) - these held cidre/jus values belonging elsewhere.
$ws.Range("B7:B10").ClearContents()

# Move the active selection to C7, matching where the user left off.
$ws.Range("C7").Select()
